# Lattice-multiplication practice sheet: refresh the random problem set.
# The document is a single 5-row x 3-column table; every cell holds one
# worked-example "card" as 5 lines (joined with manual line breaks):
#   "NN x NN"
#   "  D    D"
#   "  ----"
#   "D|    |"
#   "D|    |"
# Each of the 15 cells gets new numbers, but the card layout itself is
# unchanged, so we just overwrite each cell's Range.Text in place,
# re-joining the 5 lines with a vertical-tab (chr 11) which is how Word's
# object model represents a manual line break (<w:br/>) inside Range.Text.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$brk = [char]11

function Set-Card($row, $col, $top, $mid, $bottom1, $bottom2) {
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $top + $brk + $mid + $brk + "  ----" + $brk + $bottom1 + $brk + $bottom2
}

Set-Card 1 1 "88 x 46" "  4    6" "8|    |" "8|    |"
Set-Card 1 2 "82 x 70" "  7    0" "8|    |" "2|    |"
Set-Card 1 3 "47 x 37" "  3    7" "4|    |" "7|    |"

Set-Card 2 1 "14 x 71" "  7    1" "1|    |" "4|    |"
Set-Card 2 2 "58 x 86" "  8    6" "5|    |" "8|    |"
Set-Card 2 3 "29 x 96" "  9    6" "2|    |" "9|    |"

Set-Card 3 1 "17 x 32" "  3    2" "1|    |" "7|    |"
Set-Card 3 2 "48 x 20" "  2    0" "4|    |" "8|    |"
Set-Card 3 3 "12 x 12" "  1    2" "1|    |" "2|    |"

Set-Card 4 1 "68 x 60" "  6    0" "6|    |" "8|    |"
Set-Card 4 2 "15 x 89" "  8    9" "1|    |" "5|    |"
Set-Card 4 3 "56 x 61" "  6    1" "5|    |" "6|    |"

Set-Card 5 1 "99 x 41" "  4    1" "9|    |" "9|    |"
Set-Card 5 2 "89 x 38" "  3    8" "8|    |" "9|    |"
Set-Card 5 3 "30 x 17" "  1    7" "3|    |" "0|    |"

Write-Output "Updated 15 lattice-multiplication cards."
